# Apply the edit described by the diff:
#  - survey sheet: change the "type" column value for the FA_FOL_date,
#    FA_time_start and FA_time_end fields from "date"/"time" to "text"
#  - survey sheet: move the active selection from C12 to C11
#  - sharedStrings "date" and "time" entries become unused and are
#    dropped automatically by the engine when it re-serializes the file

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 2  -> FA_FOL_date  : date -> text
$ws.Range("C2").Value = "text"

# Row 9  -> FA_time_start : time -> text
$ws.Range("C9").Value = "text"

# Row 10 -> FA_time_end   : time -> text
$ws.Range("C10").Value = "text"

# Move the selection / active cell from C12 to C11
$ws.Range("C11").Select() | Out-Null
